$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = 'juan diego'
$ws.Range("A26").ClearFormats()
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'ramirez rendon'
$ws.Range("B26").ClearFormats()
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'ramirezrendonjuandiego54@gmail.com'
$ws.Range("C26").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '+573116347492'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '1'
$ws.Range("E26").ClearFormats()
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = '2007-08-21'
$ws.Range("F26").ClearFormats()
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '2'
$ws.Range("G26").ClearFormats()
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '[''Python'', ''JavaScript'', ''HTML'', ''CSS'', ''Flask'']'
$ws.Range("H26").ClearFormats()
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 100
$ws.Range("K26").NumberFormat = "@"
$ws.Range("K26").Value = '1'
$ws.Range("K26").ClearFormats()
$ws.Range("L26").Value = 30000000
$ws.Range("M26").NumberFormat = "@"
$ws.Range("M26").Value = 'itagui'
$ws.Range("M26").ClearFormats()
$ws.Range("N26").Value = 7
$ws.Range("O26").NumberFormat = "@"
$ws.Range("O26").Value = '[''Comunicación'', ''Adaptabilidad'', ''Creatividad'', ''Empatía'', ''Escucha activa'']'
$ws.Range("O26").ClearFormats()
$ws.Range("P26").NumberFormat = "@"
$ws.Range("P26").Value = '[''Inglés'', ''Español'']'
$ws.Range("P26").ClearFormats()
$ws.Range("Q26").NumberFormat = "@"
$ws.Range("Q26").Value = '2'
$ws.Range("Q26").ClearFormats()
$ws.Range("R26").NumberFormat = "@"
$ws.Range("R26").Value = 'backend'
$ws.Range("R26").ClearFormats()
$ws.Range("S26").NumberFormat = "@"
$ws.Range("S26").Value = '[''Autonomía'']'
$ws.Range("S26").ClearFormats()
$ws.Range("T26").Value = 345
$ws.Range("U26").NumberFormat = "@"
$ws.Range("U26").Value = 'lider_equipo'
$ws.Range("U26").ClearFormats()

# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = 'miguel'
$ws.Range("A27").ClearFormats()
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'ospina baena'
$ws.Range("B27").ClearFormats()
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'miguelospinabaena@gmail.com'
$ws.Range("C27").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '+573213983128'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '1'
$ws.Range("E27").ClearFormats()
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = '2008-02-18'
$ws.Range("F27").ClearFormats()
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '1'
$ws.Range("G27").ClearFormats()
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '[''JavaScript'', ''Java'', ''HTML'', ''CSS'']'
$ws.Range("H27").ClearFormats()
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 85
$ws.Range("K27").NumberFormat = "@"
$ws.Range("K27").Value = '3'
$ws.Range("K27").ClearFormats()
$ws.Range("L27").Value = 1623508
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = 'medellin'
$ws.Range("M27").ClearFormats()
$ws.Range("N27").Value = 3
$ws.Range("O27").NumberFormat = "@"
$ws.Range("O27").Value = '[''Comunicación'', ''Adaptabilidad'', ''Empatía'', ''Gestión del tiempo'']'
$ws.Range("O27").ClearFormats()
$ws.Range("P27").NumberFormat = "@"
$ws.Range("P27").Value = '[''Inglés'', ''Español'', ''Francés'']'
$ws.Range("P27").ClearFormats()
$ws.Range("Q27").NumberFormat = "@"
$ws.Range("Q27").Value = '2'
$ws.Range("Q27").ClearFormats()
$ws.Range("R27").NumberFormat = "@"
$ws.Range("R27").Value = 'backend'
$ws.Range("R27").ClearFormats()
$ws.Range("S27").NumberFormat = "@"
$ws.Range("S27").Value = '[''Organización'', ''Empatía'', ''Autonomía'', ''Comunicación'', ''Curiosidad'', ''Perseverancia'', ''Compromiso'']'
$ws.Range("S27").ClearFormats()
$ws.Range("T27").Value = 450
$ws.Range("U27").NumberFormat = "@"
$ws.Range("U27").Value = 'lider_equipo'
$ws.Range("U27").ClearFormats()

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = 'miguel'
$ws.Range("A28").ClearFormats()
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'ospina baena'
$ws.Range("B28").ClearFormats()
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'miguelospinabaena@gmail.com'
$ws.Range("C28").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '+573213983128'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1'
$ws.Range("E28").ClearFormats()
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = '2008-02-18'
$ws.Range("F28").ClearFormats()
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '1'
$ws.Range("G28").ClearFormats()
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '[''JavaScript'', ''Java'', ''HTML'', ''CSS'']'
$ws.Range("H28").ClearFormats()
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 85
$ws.Range("K28").NumberFormat = "@"
$ws.Range("K28").Value = '3'
$ws.Range("K28").ClearFormats()
$ws.Range("L28").Value = 1623508
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value = 'medellin'
$ws.Range("M28").ClearFormats()
$ws.Range("N28").Value = 3
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = '[''Comunicación'', ''Adaptabilidad'', ''Empatía'', ''Gestión del tiempo'']'
$ws.Range("O28").ClearFormats()
$ws.Range("P28").NumberFormat = "@"
$ws.Range("P28").Value = '[''Inglés'', ''Español'', ''Francés'']'
$ws.Range("P28").ClearFormats()
$ws.Range("Q28").NumberFormat = "@"
$ws.Range("Q28").Value = '2'
$ws.Range("Q28").ClearFormats()
$ws.Range("R28").NumberFormat = "@"
$ws.Range("R28").Value = 'backend'
$ws.Range("R28").ClearFormats()
$ws.Range("S28").NumberFormat = "@"
$ws.Range("S28").Value = '[''Organización'', ''Empatía'', ''Autonomía'', ''Comunicación'', ''Curiosidad'', ''Perseverancia'', ''Compromiso'']'
$ws.Range("S28").ClearFormats()
$ws.Range("T28").Value = 450
$ws.Range("U28").NumberFormat = "@"
$ws.Range("U28").Value = 'lider_equipo'
$ws.Range("U28").ClearFormats()

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = 'miguel'
$ws.Range("A29").ClearFormats()
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'ospina baena'
$ws.Range("B29").ClearFormats()
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'miguelospinabaena@gmail.com'
$ws.Range("C29").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '+573213983128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '1'
$ws.Range("E29").ClearFormats()
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = '2008-02-18'
$ws.Range("F29").ClearFormats()
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '1'
$ws.Range("G29").ClearFormats()
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = '[''JavaScript'', ''Java'', ''HTML'', ''CSS'']'
$ws.Range("H29").ClearFormats()
$ws.Range("I29").Value = 8
$ws.Range("J29").Value = 85
$ws.Range("K29").NumberFormat = "@"
$ws.Range("K29").Value = '3'
$ws.Range("K29").ClearFormats()
$ws.Range("L29").Value = 1623508
$ws.Range("M29").NumberFormat = "@"
$ws.Range("M29").Value = 'medellin'
$ws.Range("M29").ClearFormats()
$ws.Range("N29").Value = 3
$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = '[''Comunicación'', ''Adaptabilidad'', ''Empatía'', ''Gestión del tiempo'']'
$ws.Range("O29").ClearFormats()
$ws.Range("P29").NumberFormat = "@"
$ws.Range("P29").Value = '[''Inglés'', ''Español'', ''Francés'']'
$ws.Range("P29").ClearFormats()
$ws.Range("Q29").NumberFormat = "@"
$ws.Range("Q29").Value = '2'
$ws.Range("Q29").ClearFormats()
$ws.Range("R29").NumberFormat = "@"
$ws.Range("R29").Value = 'backend'
$ws.Range("R29").ClearFormats()
$ws.Range("S29").NumberFormat = "@"
$ws.Range("S29").Value = '[''Organización'', ''Empatía'', ''Autonomía'', ''Comunicación'', ''Curiosidad'', ''Perseverancia'', ''Compromiso'']'
$ws.Range("S29").ClearFormats()
$ws.Range("T29").Value = 450
$ws.Range("U29").NumberFormat = "@"
$ws.Range("U29").Value = 'lider_equipo'
$ws.Range("U29").ClearFormats()

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = 'miguel'
$ws.Range("A30").ClearFormats()
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'ospina baena'
$ws.Range("B30").ClearFormats()
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'miguelospinabaena@gmail.com'
$ws.Range("C30").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '+573213983128'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '1'
$ws.Range("E30").ClearFormats()
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = '2008-02-18'
$ws.Range("F30").ClearFormats()
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '1'
$ws.Range("G30").ClearFormats()
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '[''JavaScript'', ''Java'', ''HTML'', ''CSS'']'
$ws.Range("H30").ClearFormats()
$ws.Range("I30").Value = 8
$ws.Range("J30").Value = 85
$ws.Range("K30").NumberFormat = "@"
$ws.Range("K30").Value = '3'
$ws.Range("K30").ClearFormats()
$ws.Range("L30").Value = 1623508
$ws.Range("M30").NumberFormat = "@"
$ws.Range("M30").Value = 'medellin'
$ws.Range("M30").ClearFormats()
$ws.Range("N30").Value = 3
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = '[''Comunicación'', ''Adaptabilidad'', ''Empatía'', ''Gestión del tiempo'']'
$ws.Range("O30").ClearFormats()
$ws.Range("P30").NumberFormat = "@"
$ws.Range("P30").Value = '[''Inglés'', ''Español'', ''Francés'']'
$ws.Range("P30").ClearFormats()
$ws.Range("Q30").NumberFormat = "@"
$ws.Range("Q30").Value = '2'
$ws.Range("Q30").ClearFormats()
$ws.Range("R30").NumberFormat = "@"
$ws.Range("R30").Value = 'backend'
$ws.Range("R30").ClearFormats()
$ws.Range("S30").NumberFormat = "@"
$ws.Range("S30").Value = '[''Organización'', ''Empatía'', ''Autonomía'', ''Comunicación'', ''Curiosidad'', ''Perseverancia'', ''Compromiso'']'
$ws.Range("S30").ClearFormats()
$ws.Range("T30").Value = 450
$ws.Range("U30").NumberFormat = "@"
$ws.Range("U30").Value = 'lider_equipo'
$ws.Range("U30").ClearFormats()

# Row 31
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = 'juan diego'
$ws.Range("A31").ClearFormats()
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'ramirez rendon'
$ws.Range("B31").ClearFormats()
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'ramirezrendonjuandiego5@gmail.com'
$ws.Range("C31").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '+573116347491'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '1'
$ws.Range("E31").ClearFormats()
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = '2007-08-21'
$ws.Range("F31").ClearFormats()
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '1'
$ws.Range("G31").ClearFormats()
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '[''AWS'']'
$ws.Range("H31").ClearFormats()
$ws.Range("I31").Value = 10
$ws.Range("J31").Value = 100
$ws.Range("K31").NumberFormat = "@"
$ws.Range("K31").Value = '2'
$ws.Range("K31").ClearFormats()
$ws.Range("L31").Value = 1623500
$ws.Range("M31").NumberFormat = "@"
$ws.Range("M31").Value = 'itagui'
$ws.Range("M31").ClearFormats()
$ws.Range("N31").Value = 7
$ws.Range("O31").NumberFormat = "@"
$ws.Range("O31").Value = '[''Comunicación'']'
$ws.Range("O31").ClearFormats()
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = '[''Español'']'
$ws.Range("P31").ClearFormats()
$ws.Range("Q31").NumberFormat = "@"
$ws.Range("Q31").Value = '3'
$ws.Range("Q31").ClearFormats()
$ws.Range("R31").NumberFormat = "@"
$ws.Range("R31").Value = 'frontend'
$ws.Range("R31").ClearFormats()
$ws.Range("S31").NumberFormat = "@"
$ws.Range("S31").Value = '[''Autonomía'']'
$ws.Range("S31").ClearFormats()
$ws.Range("T31").Value = 345
$ws.Range("U31").NumberFormat = "@"
$ws.Range("U31").Value = 'emprender'
$ws.Range("U31").ClearFormats()

Write-Output "Rows 26-31 added"